# Applies the "break out stock.yaml completed" edit to the weekly
# COALINDIA.NS stock-history sheet:
#   1. Resets a batch of previously-computed "backup" (R) flags (and one
#      "detect_structure" / Q flag) back to 0 now that the structure
#      detection has been re-run / broken out.
#   2. Fixes up the most recent trailing rows (711/713/714) whose
#      isPivot/backup values were placeholders.
#   3. Appends 13 new weekly OHLCV rows (715-727) that were pulled in by
#      the refreshed download, extending the sheet through 2024-09-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rows whose "backup" (column R) flag flips from 1/2 back to 0.
# ---------------------------------------------------------------------
$rFlagRows = @(56, 67, 71, 85, 88, 99, 106, 108, 114, 117, 143, 155, 160, `
    168, 175, 196, 209, 214, 217, 228, 230, 237, 264, 273, 280, 282, 292, `
    329, 338, 359, 361, 372, 379, 390, 393, 401, 403, 416, 425, 440, 453, `
    456, 467, 478, 483, 488, 498, 506, 508, 514, 520, 527, 542, 546, 554, `
    561, 572, 578, 589, 598, 606, 609, 626, 629, 635, 646, 664, 701, 705)

foreach ($r in $rFlagRows) {
    $ws.Cells.Item($r, 18).Value = 0   # column R = backup
}

# Row 56 also had its "detect_structure" (column Q) flag reset.
$ws.Cells.Item(56, 17).Value = 0       # column Q = detect_structure

# ---------------------------------------------------------------------
# 2. Trailing-row fixups.
# ---------------------------------------------------------------------
$ws.Cells.Item(711, 15).Value = 3      # column O = isPivot
$ws.Cells.Item(713, 18).Value = 0      # column R = backup (was blank)
$ws.Cells.Item(714, 18).Value = 0      # column R = backup (was blank)

# ---------------------------------------------------------------------
# 3. Append the 13 new weekly rows (715-727). Columns A-Q are filled in;
#    column R ("backup") is left blank, matching the newly pulled rows.
# ---------------------------------------------------------------------
$newRows = @(
    @(715, 45474, 473.3500061035156, 493.2999877929688, 471.25, 491.5, 486.6307678222656, 44005271, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(716, 45481, 492.7000122070312, 508.6000061035156, 480.5499877929688, 496.2000122070312, 491.2842102050781, 48248769, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(717, 45488, 500, 518.4000244140625, 486.75, 488, 483.1654357910156, 52788345, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(718, 45495, 488, 511.7000122070312, 464.5499877929688, 509.8500061035156, 504.7989807128906, 56135857, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(719, 45502, 514.5, 542.25, 507.5, 524.5, 519.3038330078125, 78178883, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(720, 45509, 515, 535, 497.5499877929688, 529.7999877929688, 524.5513305664062, 57363442, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(721, 45516, 529, 530.5, 499.0499877929688, 512.2999877929688, 507.2247009277344, 32520348, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(722, 45523, 513.7999877929688, 541.1500244140625, 513.1500244140625, 538.8499755859375, 538.8499755859375, 49094981, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(723, 45530, 542, 543.5499877929688, 516.0999755859375, 524.9500122070312, 524.9500122070312, 61194214, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(724, 45537, 527, 529, 481.0499877929688, 488.5499877929688, 488.5499877929688, 52620510, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(725, 45544, 488, 497.8500061035156, 478.0499877929688, 490.25, 490.25, 48286152, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(726, 45551, 492.5, 498, 476.25, 490.9500122070312, 490.9500122070312, 39832969, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(727, 45558, 493.5, 517.8499755859375, 490.5, 516.0999755859375, 516.0999755859375, 40809801, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A (Datetime) keeps the same date number-format style as the
    # rest of the column.
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[2]    # Open
    $ws.Cells.Item($r, 3).Value = $row[3]    # High
    $ws.Cells.Item($r, 4).Value = $row[4]    # Low
    $ws.Cells.Item($r, 5).Value = $row[5]    # Close
    $ws.Cells.Item($r, 6).Value = $row[6]    # Adj Close
    $ws.Cells.Item($r, 7).Value = $row[7]    # Volume
    $ws.Cells.Item($r, 8).Value = $row[8]    # Year
    $ws.Cells.Item($r, 9).Value = $row[9]    # Month
    $ws.Cells.Item($r, 10).Value = $row[10]  # Day
    $ws.Cells.Item($r, 11).Value = $row[11]  # Hour
    $ws.Cells.Item($r, 12).Value = $row[12]  # Minute
    $ws.Cells.Item($r, 13).Value = $row[13]  # Second
    $ws.Cells.Item($r, 14).Value = $row[14]  # Week
    $ws.Cells.Item($r, 15).Value = $row[15]  # isPivot
    $ws.Cells.Item($r, 16).Value = $row[16]  # two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[17]  # detect_structure
    # column 18 (R / backup) intentionally left blank for the new rows.
}
